$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 47 (shifts existing rows 47:82 down to 48:83,
# dimension grows to A1:R83 automatically).
$ws.Rows("47").Insert()

# Populate the new row 47 with this week's record (same market/category
# as its neighbours; only the date + price/volume figures are new).
$ws.Range("A47").Value = 11
$ws.Range("B47").Value = "Vega Monumental Concepción"
$ws.Range("C47").Value = "Bíobío"
$ws.Range("D47").Value = 44664
$ws.Range("E47").Value = 8
$ws.Range("F47").Value = 100112001
$ws.Range("G47").Value = "Berenjena"
$ws.Range("H47").Value = "Sin especificar"
$ws.Range("I47").Value = "Primera"
$ws.Range("J47").Value = 100
$ws.Range("K47").Value = 6500
$ws.Range("L47").Value = 7000
$ws.Range("M47").Value = 6750
$ws.Range("N47").Value = "$/caja 60 unidades"
$ws.Range("O47").Value = "Región de Arica y Parinacota"
$ws.Range("P47").Value = 112
$ws.Range("Q47").Value = 60
$ws.Range("R47").Value = "Hortaliza"
